$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows related to "Categoria de Ingresso" (ticket category) use cases:
# CSU-04 (row 5), CSU-14 (row 15), CSU-15 (row 16), CSU-16 (row 17)
# Delete from bottom to top so the remaining row indices stay valid.
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(5).Delete()

# Update selection to match the recorded view state after the edit.
$ws.Range("B14").Select()
